# "Reviewed CYRS & HSI" - update the PO3 DGW PP schedule:
#  - HSI (row 7)               -> mark Done, fill Actual End Date
#  - CYRS Document Update (r9) -> mark Done, fill Actual End Date
#  - CYRS & HSI Review (r10)   -> mark Done, fill Actual/Expected dates, clear comment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: HSI ---
$c = $ws.Range("F7")
$c.Value = 43923
$c.NumberFormat = "mm-dd-yy"
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108
$ws.Range("G7").Value = "Done"

# --- Row 9: CYRS Document Update ---
$c = $ws.Range("F9")
$c.Value = 43923
$c.NumberFormat = "mm-dd-yy"
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108
$ws.Range("G9").Value = "Done"

# --- Row 10: CYRS & HSI Review ---
$c = $ws.Range("E10")
$c.Value = 43953
$c.NumberFormat = "mm-dd-yy"
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

$c = $ws.Range("F10")
$c.Value = 43953
$c.NumberFormat = "mm-dd-yy"
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

$ws.Range("G10").Value = "Done"
$ws.Range("H10").ClearContents()

# --- Selection, as left by the author after editing E11 ---
$ws.Range("E11").Select()
